$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "group by month" logic is buggy: the "fery" transaction (row 3) gets dropped
# from the report entirely, leaving that row blank.
$ws.Range("A3:E3").ClearContents()

# The Total row moves up from row 5 to row 4. Because of the grouping bug the
# "Uang Masuk" (money in) total no longer includes the dropped row, so it
# becomes 0, while "Uang Keluar" (money out) stays at 5000.
$ws.Range("C4").Value = "Total"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 5000

# Old Total row (row 5) is now empty.
$ws.Range("A5:E5").ClearContents()
